$d = $word.ActiveDocument

# Locate the paragraph that begins "Unless defendant was previously fingerprinted ..."
# and prepend a new bold run reading "Fingerprinting Required: " right before it,
# matching the paragraph's existing run formatting (Palatino Linotype, 10pt).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Unless defendant was previously fingerprinted", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPoint = $d.Range($searchRange.Start, $searchRange.Start)

    $openXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r>' +
        '<w:rPr>' +
        '<w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">Fingerprinting Required: </w:t>' +
        '</w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $insertPoint.InsertXML($openXml)
}
